$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.592835187911987
$ws.Range("B1").Value = 3.39194917678833
$ws.Range("C1").Value = 5.57266902923584
$ws.Range("D1").Value = 6.088551044464111
$ws.Range("E1").Value = 1.099314212799072
